# edit.ps1 - applies the documented change:
#   1. The old "_GoBack" bookmark (sitting right after "14.03", before
#      ".2016") is removed - Word keeps only a single "_GoBack" bookmark
#      in the whole document, tracking the most recent edit location.
#   2. "postoji id." becomes "postoji itd." (a "t" is typed in the
#      middle of "id.", turning the typo into the abbreviation "itd."),
#      and the new "_GoBack" bookmark ends up right after the freshly
#      typed "t" (i.e. directly before "d."), matching where Word
#      would drop it after that keystroke.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: remove the old "_GoBack" bookmark located right after "14.03"
# ---------------------------------------------------------------------
$oldMarker = $d.Content
$oldMarker.Find.Execute("14.03", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0)

# The bookmark sits exactly at the end of the found "14.03" text. Lift
# one character spanning that boundary out and drop it back in
# unchanged - that forces the (now orphaned) bookmark pair to be
# dropped while leaving the visible text identical.
$bmSpan = $d.Range($oldMarker.End - 1, $oldMarker.End + 1)
$bmSpanText = $bmSpan.Text
$bmSpan.Delete()
$d.Range($oldMarker.End - 1, $oldMarker.End - 1).InsertBefore($bmSpanText)

# ---------------------------------------------------------------------
# Step 2: "postoji id." -> "postoji itd.", with a fresh "_GoBack"
#         bookmark right after the newly-typed "t"
# ---------------------------------------------------------------------
$target = $d.Content
$target.Find.Execute("postoji id.", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0)

# Position right between the "i" and the "d" of "id."
$caret = $d.Range($target.Start + 9, $target.Start + 9)
$caret.InsertBefore("t")

# Force the newly typed "t" onto its own run (mirrors the run split
# Word performs as you type), then park the "_GoBack" bookmark right
# after it, i.e. right before "d."
$typedT = $d.Range($target.Start + 9, $target.Start + 10)
$typedT.Bold = 1
$typedT.Bold = 0

$newMarker = $d.Range($target.Start + 10, $target.Start + 10)
$d.Bookmarks.Add("_GoBack", $newMarker)
